$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSE")

# ---------------------------------------------------------------------------
# 1) Copy the D3:F57 "Random Forest-100 (superdataset-20.csv)" block formats
#    over to I3:K57 so the new columns inherit the same header / number
#    styles (bold-red header, bold-centered row numbers, scientific data).
#    Done as two pieces (3:54 and 56:57) so the row-3-to-54 gap at row 55
#    (no row 55 exists in the sheet) doesn't get materialised as an empty
#    row in the destination.
# ---------------------------------------------------------------------------
$ws.Range("D3:F54").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D56:F57").Copy()
$ws.Range("I56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Paste-formats of row 3 (D3:F3) also stamps a blank placeholder on J3
# (mirroring the always-empty E3) - the target sheet has no J3 cell at all,
# so drop it.
$ws.Range("J3").Clear()

# ---------------------------------------------------------------------------
# 2) New column widths for J/K (mirrors E/F widths for the new block).
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 13.140625
$ws.Columns.Item(11).ColumnWidth = 11.9

# ---------------------------------------------------------------------------
# 3) D3 now reads "...without cons" (shared string already exists), and the
#    previous "Random Forest-100 (superdataset-20.csv)" header moves to I3.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "Random Forest-100 (superdataset-20.csv without cons)"
$ws.Range("I3").Value = "Random Forest-100 (superdataset-20.csv)"

# ---------------------------------------------------------------------------
# 4) Row 4 headers for the new block (train (MSE) / test (MSE)).
# ---------------------------------------------------------------------------
$ws.Range("J4").Value = "train (MSE)"
$ws.Range("K4").Value = "test (MSE)"

# ---------------------------------------------------------------------------
# 5) Row numbers 1..50 in column I (D6:D54 pattern mirrored), J/K stay blank
#    just like E/F used to be before this edit.
# ---------------------------------------------------------------------------
$ws.Range("I5").Value = 1
$ws.Range("I6").Formula = "=I5+1"
$ws.Range("I7:I54").Formula = "=I6+1"

# ---------------------------------------------------------------------------
# 6) avg / SD rows for the new (empty) I/J/K block - same formulas as D/E/F,
#    which will legitimately evaluate to #DIV/0! since J5:J54/K5:K54 are
#    still empty.
# ---------------------------------------------------------------------------
$ws.Range("I56").Value = "avg"
$ws.Range("J56").Formula = "=AVERAGE(J5:J54)"
$ws.Range("K56").Formula = "=AVERAGE(K5:K54)"
$ws.Range("I57").Value = "SD"
$ws.Range("J57").Formula = "=_xlfn.STDEV.S(J5:J54)"
$ws.Range("K57").Formula = "=_xlfn.STDEV.S(K5:K54)"

# ---------------------------------------------------------------------------
# 7) Fill in the actual train/test MSE values for the original D/E/F block
#    (previously blank placeholders).
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = 0.0001244664489460285
$ws.Range("F5").Value = 0.0006179644525098096
$ws.Range("E6").Value = 0.0001168423218238618
$ws.Range("F6").Value = 0.0009055531980677242
$ws.Range("E7").Value = 0.000130648119186104
$ws.Range("F7").Value = 0.0006158128562131259
$ws.Range("E8").Value = 0.0001290726638349093
$ws.Range("F8").Value = 0.0008222867701021917
$ws.Range("E9").Value = 0.0001196053616554614
$ws.Range("F9").Value = 0.001131388841680097
$ws.Range("E10").Value = 0.000115076839691538
$ws.Range("F10").Value = 0.0008787776240302438
$ws.Range("E11").Value = 0.0001317934513385996
$ws.Range("F11").Value = 0.0008121420180910622
$ws.Range("E12").Value = 0.000127837345502081
$ws.Range("F12").Value = 0.0006582387370480782
$ws.Range("E13").Value = 0.0001351005404432482
$ws.Range("F13").Value = 0.0007972935365645491
$ws.Range("E14").Value = 0.0001278171048235324
$ws.Range("F14").Value = 0.0008619840199135268
$ws.Range("E15").Value = 0.0001244327257755087
$ws.Range("F15").Value = 0.0009167811602896157
$ws.Range("E16").Value = 0.000121383692869805
$ws.Range("F16").Value = 0.0007402899750111889
$ws.Range("E17").Value = 0.000133103359281455
$ws.Range("F17").Value = 0.0007114286833073213
$ws.Range("E18").Value = 0.0001258810634120846
$ws.Range("F18").Value = 0.000966483749525564
$ws.Range("E19").Value = 0.000131473362263601
$ws.Range("F19").Value = 0.0006312700112743109
$ws.Range("E20").Value = 0.0001257135650161644
$ws.Range("F20").Value = 0.0008452603224812575
$ws.Range("E21").Value = 0.0001220539330548962
$ws.Range("F21").Value = 0.0011109788447526962
$ws.Range("E22").Value = 0.0001102856884126588
$ws.Range("F22").Value = 0.0011593087544397052
$ws.Range("E23").Value = 0.0001227221323945123
$ws.Range("F23").Value = 0.0006402535770419681
$ws.Range("E24").Value = 0.0001314882059188829
$ws.Range("F24").Value = 0.0006921951238210829
$ws.Range("E25").Value = 0.000117294958719843
$ws.Range("F25").Value = 0.001011404424568403
$ws.Range("E26").Value = 0.0001244882008659378
$ws.Range("F26").Value = 0.0006885361642005277
$ws.Range("E27").Value = 0.0001225717054711288
$ws.Range("F27").Value = 0.0008471208776922449
$ws.Range("E28").Value = 0.0001191280562107632
$ws.Range("F28").Value = 0.0009263166421307512
$ws.Range("E29").Value = 0.0001270736768180333
$ws.Range("F29").Value = 0.001052457409967145
$ws.Range("E30").Value = 0.000110120986333765
$ws.Range("F30").Value = 0.001066169818548743
$ws.Range("E31").Value = 0.0001163141197106293
$ws.Range("F31").Value = 0.0008225201401958283
$ws.Range("E32").Value = 0.000116918167073448
$ws.Range("F32").Value = 0.0008418710725765221
$ws.Range("E33").Value = 0.000118660276321481
$ws.Range("F33").Value = 0.000843904520480517
$ws.Range("E34").Value = 0.0001158110337037697
$ws.Range("F34").Value = 0.000901476902262500
$ws.Range("E35").Value = 0.0001165808712408225
$ws.Range("F35").Value = 0.001025380525361276
$ws.Range("E36").Value = 0.000114913852455661
$ws.Range("F36").Value = 0.001093108640349606
$ws.Range("E37").Value = 0.000116172173264466
$ws.Range("F37").Value = 0.0009031327797790908
$ws.Range("E38").Value = 0.0001201380420118854
$ws.Range("F38").Value = 0.001109221952075963
$ws.Range("E39").Value = 0.0001093439754658021
$ws.Range("F39").Value = 0.001247069665065622
$ws.Range("E40").Value = 0.000132007132184734
$ws.Range("F40").Value = 0.0006064184238408613
$ws.Range("E41").Value = 0.000124389307875767
$ws.Range("F41").Value = 0.0006413728466023837
$ws.Range("E42").Value = 0.000127969150138329
$ws.Range("F42").Value = 0.0007842184876207509
$ws.Range("E43").Value = 0.00012662320151519
$ws.Range("F43").Value = 0.0008467408967656146
$ws.Range("E44").Value = 0.000114893492568111
$ws.Range("F44").Value = 0.001130988879419156
$ws.Range("E45").Value = 0.0001250342538163027
$ws.Range("F45").Value = 0.0008125776212513974
$ws.Range("E46").Value = 0.0001029452361554765
$ws.Range("F46").Value = 0.00146714595666021
$ws.Range("E47").Value = 0.000116840757768657
$ws.Range("F47").Value = 0.0008789791007115342
$ws.Range("E48").Value = 0.000122162218744710
$ws.Range("F48").Value = 0.0009862295153737317
$ws.Range("E49").Value = 0.000118763128036097
$ws.Range("F49").Value = 0.001057709492025047
$ws.Range("E50").Value = 0.000131227511248206
$ws.Range("F50").Value = 0.0008421149844078173
$ws.Range("E51").Value = 0.0001247137236947421
$ws.Range("F51").Value = 0.0007455592373415175
$ws.Range("E52").Value = 0.0001113379617603689
$ws.Range("F52").Value = 0.0009565909473930237
$ws.Range("E53").Value = 0.000116765444734906
$ws.Range("F53").Value = 0.0009518850530888231
$ws.Range("E54").Value = 0.000127919134939020
$ws.Range("F54").Value = 0.0008439628751452543

# ---------------------------------------------------------------------------
# 8) Restore the view: scroll so row 22 is at the top and M40 is selected
#    (matches the author re-positioning the viewport on the MSE sheet).
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A22"))
$ws.Range("M40").Select()
